# "finished commenting compilation, ensuring it knits still"
#
# The working-hours log on Sheet1 gets one more day entered (Sunday
# 2024-07-14, previously missing) ahead of the existing last row
# (Monday 2024-07-15), and that Monday's hours/notes get updated to
# reflect the extra work actually done that day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Make room for the new Sunday row just above the current last row
# (currently row 51, the Monday 2024-07-15 entry) so it becomes row 52.
$ws.Rows.Item(51).Insert()

# New row 51: Sunday 2024-07-14 -- 4 hours.
$ws.Range("A51").Value = 45487
$ws.Range("A51").NumberFormat = "d-mmm"
$ws.Range("B51").Value = "S"
$ws.Range("B51").NumberFormat = "d-mmm"
$ws.Range("C51").Value = 4
$ws.Range("E51").Value = "Differentiated the three studentgrades (_rep and _prof). Started commenting code, some revising on code."
$ws.Range("E51").WrapText = $true

# Row 52 (the pre-existing Monday 2024-07-15 row, now shifted down):
# hours revised from 2 to 8 and notes on finishing up the commenting added.
$ws.Range("C52").Value = 8
$ws.Range("E52").Value = "Finished commenting code. Ensuring it runs/knits correctly. Need to finish up error bar section, RFE, and xgbms"
$ws.Range("E52").WrapText = $true

# Leave the selection on the note just typed, as the author did.
$null = $ws.Range("E52").Select()
